$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab title to "Inventory"
$ws.Name = "Inventory"

# Update the Stock column (C) values
$ws.Range("C2").Value = 249
$ws.Range("C3").Value = 787
$ws.Range("C4").Value = 1537
$ws.Range("C5").Value = 1251
$ws.Range("C6").Value = 1210
$ws.Range("C7").Value = 1382
$ws.Range("C8").Value = 472
$ws.Range("C9").Value = 594
$ws.Range("C10").Value = 1255
$ws.Range("C11").Value = 610
$ws.Range("C12").Value = 1187
$ws.Range("C13").Value = 679
$ws.Range("C14").Value = 1387
$ws.Range("C15").Value = 807
$ws.Range("C16").Value = 537
$ws.Range("C17").Value = 1681
$ws.Range("C18").Value = 871
$ws.Range("C19").Value = 1622
$ws.Range("C20").Value = 645
$ws.Range("C21").Value = 1565
$ws.Range("C22").Value = 755
$ws.Range("C23").Value = 1476
$ws.Range("C24").Value = 400
$ws.Range("C25").Value = 1795
